# "updates with supplemental figures"
# Refreshes the Fixed-Residue PSSM grid (Sheet1!B2:K21) with the recomputed
# log-odds scores. Layout/headers/row-labels (row 1, column A) are unchanged;
# only the numeric matrix body is rewritten cell-by-cell to match the new
# PSSM values (including the new "no data" floor of -17.3265874839145,
# replacing the previous -19.7552934769669).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values per row (columns B..K = column indices 2..11)
$rowValues = @{
    2  = @(-17.3265874839145, 2.562169140382605, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    3  = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, 2.796115411921812, -17.3265874839145, -17.3265874839145)
    4  = @(-17.3265874839145, 2.412535497788493, 1.930093442580538, -17.3265874839145, 3.149922403991999, -17.3265874839145, 1.63009778591346, -17.3265874839145, 2.372010368362444, -17.3265874839145)
    5  = @(-17.3265874839145, 1.299624188840759, -17.3265874839145, -17.3265874839145, -17.3265874839145, 2.915671366110815, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    6  = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    7  = @(2.81791268665524, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    8  = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    9  = @(3.694710401472615, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    10 = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, 1.247787996701922, -17.3265874839145, 1.865033245754724)
    11 = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, 2.676552938465519, -17.3265874839145, -17.3265874839145, -17.3265874839145, 1.688067797348656)
    12 = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    13 = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, 2.002019517493798, 1.929540563167561)
    14 = @(-17.3265874839145, -17.3265874839145, 1.512951566806012, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, 2.197083579792345)
    15 = @(-17.3265874839145, -17.3265874839145, 1.979894805610651, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145)
    16 = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, 2.205937763628035, -17.3265874839145)
    17 = @(-17.3265874839145, 1.572890405993413, 1.737073548690091, -17.3265874839145, -17.3265874839145, -17.3265874839145, 1.066136612321032, 1.738126339928046, 1.802245123700718, -17.3265874839145)
    18 = @(-17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, -17.3265874839145, 1.369975993110649, 1.025839641385908, 1.441651773284242, -17.3265874839145)
    19 = @(-17.3265874839145, -17.3265874839145, 1.713287406934385, -17.3265874839145, -17.3265874839145, -17.3265874839145, 1.591416778370496, 1.648580366913193, -17.3265874839145, -17.3265874839145)
    20 = @(-17.3265874839145, 0.728679061853971, 1.473987009728952, -17.3265874839145, 3.475579379373403, -17.3265874839145, 2.147389712483058, 1.119061926694944, -17.3265874839145, 2.244900391394505)
    21 = @(-17.3265874839145, 0.7451105459518594, -17.3265874839145, 4.321919756600269, -17.3265874839145, 2.599518813374347, 2.257730970025789, -17.3265874839145, -17.3265874839145, -17.3265874839145)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item([int]$r, $i + 2).Value = $vals[$i]
    }
}
